$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RQ8")
$ws.Columns("W:W").Delete()
$ws.Columns("S:S").Delete()
